$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 5.133420467376709
$ws.Range("C2").Value = 0.01023483276367188
$ws.Range("D2").Value = 0.6066503524780273
$ws.Range("E2").Value = 5.788755893707275

# Row 3
$ws.Range("B3").Value = 0.2404730319976807
$ws.Range("E3").Value = 0.2404730319976807

# Row 4
$ws.Range("B4").Value = 52.86789798736572
$ws.Range("C4").Value = 0.07901811599731445
$ws.Range("D4").Value = 5.585498094558716
$ws.Range("E4").Value = 58.88949418067932
$ws.Range("F4").Value = 1564

# Row 5
$ws.Range("B5").Value = 33.35113787651062
$ws.Range("C5").Value = 0.02010226249694824
$ws.Range("D5").Value = 1.087490320205688
$ws.Range("E5").Value = 34.52986168861389

# Row 6
$ws.Range("B6").Value = 47.25594401359558
$ws.Range("C6").Value = 0.07221770286560059
$ws.Range("D6").Value = 4.555901527404785
$ws.Range("E6").Value = 52.17206954956055

# Row 7
$ws.Range("B7").Value = 0.1918485164642334
$ws.Range("C7").Value = 0.00099945068359375
$ws.Range("E7").Value = 0.1928479671478271

# Row 8
$ws.Range("B8").Value = 0.199887752532959
$ws.Range("E8").Value = 0.199887752532959

$wb.Save()
